$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 8009.84
$ws.Range("I15").Value = 8009.84
$ws.Range("K15").Value = 24029.52
$ws.Range("M15").Value = -23860.52

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5671.5713
$ws.Range("J40").Value = 7714.2856
$ws.Range("L40").Value = 7714.2856
$ws.Range("N40").Value = -8064.2856

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 73.5
$ws.Range("I58").Value = 73.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 220.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -70.5
$ws.Range("N58").Value = ""

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").Value = ""

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 929.875
$ws.Range("I107").Value = 1109.9231
$ws.Range("K107").Value = 1109.9231
$ws.Range("M107").Value = 810.0769

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 166670980
$ws.Range("I113").Value = 200002370
$ws.Range("K113").Value = 200002370
$ws.Range("M113").Value = -199999116

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 23258504
$ws.Range("J138").Value = 29414668
$ws.Range("L138").Value = 88244004
$ws.Range("N138").Value = -88254284

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5923.6665
$ws.Range("I32").Value = 5734.07
$ws.Range("K32").Value = 5734.07
$ws.Range("M32").Value = -5447.07

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 27785938
$ws.Range("I61").Value = 33340126
$ws.Range("K61").Value = 33340126
$ws.Range("M61").Value = -33339914

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3909.6572
$ws.Range("I74").Value = 2953.52
$ws.Range("K74").Value = 2953.52
$ws.Range("M74").Value = -2079.52

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3909.6572
$ws.Range("I77").Value = 2953.52
$ws.Range("K77").Value = 14767.6
$ws.Range("M77").Value = -10399.6

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3644.658
$ws.Range("I132").Value = 3013.8057
$ws.Range("K132").Value = 9041.417099999999
$ws.Range("M132").Value = -6511.417099999999

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 27785938
$ws.Range("I136").Value = 33340126
$ws.Range("K136").Value = 100020378
$ws.Range("M136").Value = -100017828

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2827.6667
$ws.Range("I22").Value = 4962.5
$ws.Range("J22").Value = 387.85715
$ws.Range("K22").Value = 4962.5
$ws.Range("L22").Value = 387.85715
$ws.Range("M22").Value = -4789.5
$ws.Range("N22").Value = -733.85715

# BSM row 70
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 149900
$ws.Range("J70").Value = 149900
$ws.Range("L70").Value = 149900
$ws.Range("N70").Value = -150486

# BSM row 73
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H73").Value = 149900
$ws.Range("J73").Value = 149900
$ws.Range("L73").Value = 149900
$ws.Range("N73").Value = -151928

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1247.0435
$ws.Range("I86").Value = 1349.1
$ws.Range("J86").Value = 566.6667
$ws.Range("K86").Value = 1349.1
$ws.Range("L86").Value = 566.6667
$ws.Range("M86").Value = -226.0999999999999
$ws.Range("N86").Value = -2812.6667

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1247.0435
$ws.Range("I89").Value = 1349.1
$ws.Range("J89").Value = 566.6667
$ws.Range("K89").Value = 6745.5
$ws.Range("L89").Value = 2833.3335
$ws.Range("M89").Value = -1129.5
$ws.Range("N89").Value = -14065.3335

# BSM row 100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 18083.5
$ws.Range("J100").Value = 18083.5
$ws.Range("L100").Value = 18083.5
$ws.Range("N100").Value = -20247.5

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4425.826
$ws.Range("I107").Value = 4053.25
$ws.Range("K107").Value = 4053.25
$ws.Range("M107").Value = -2133.25

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5958.6665
$ws.Range("I134").Value = 5670
$ws.Range("K134").Value = 17010
$ws.Range("M134").Value = -14475

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3086.2083
$ws.Range("J16").Value = 3783.4546
$ws.Range("L16").Value = 3783.4546
$ws.Range("N16").Value = -4357.4546

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 817.1818
$ws.Range("I22").Value = 420
$ws.Range("J22").Value = 1512.25
$ws.Range("K22").Value = 420
$ws.Range("L22").Value = 1512.25
$ws.Range("M22").Value = -70
$ws.Range("N22").Value = -2212.25

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8860.909
$ws.Range("I31").Value = 7500
$ws.Range("J31").Value = 10494
$ws.Range("K31").Value = 7500
$ws.Range("L31").Value = 10494
$ws.Range("M31").Value = -7205
$ws.Range("N31").Value = -11084

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8860.909
$ws.Range("I34").Value = 7500
$ws.Range("J34").Value = 10494
$ws.Range("K34").Value = 7500
$ws.Range("L34").Value = 10494
$ws.Range("M34").Value = -7298
$ws.Range("N34").Value = -10898

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9346
$ws.Range("I58").Value = 4499.8335
$ws.Range("J58").Value = 13499.857
$ws.Range("K58").Value = 4499.8335
$ws.Range("L58").Value = 13499.857
$ws.Range("M58").Value = -4296.8335
$ws.Range("N58").Value = -13905.857

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 575.25
$ws.Range("I107").Value = 518.4
$ws.Range("J107").Value = 632.1
$ws.Range("K107").Value = 518.4
$ws.Range("L107").Value = 632.1
$ws.Range("M107").Value = 1401.6
$ws.Range("N107").Value = -4472.1

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3086.2083
$ws.Range("J113").Value = 3783.4546
$ws.Range("L113").Value = 3783.4546
$ws.Range("N113").Value = -8123.4546

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 9346
$ws.Range("I136").Value = 4499.8335
$ws.Range("J136").Value = 13499.857
$ws.Range("K136").Value = 13499.5005
$ws.Range("L136").Value = 40499.571
$ws.Range("M136").Value = -10949.5005
$ws.Range("N136").Value = -45599.571

# CUL row 9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 414.14285
$ws.Range("J9").Value = 414.14285
$ws.Range("L9").Value = 1242.42855
$ws.Range("N9").Value = -1690.42855

# CUL row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 877.25
$ws.Range("I22").Value = 877.25
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2631.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2462.75
$ws.Range("N22").Value = ""

# CUL row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 877.25
$ws.Range("I27").Value = 877.25
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 2631.75
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -2529.75
$ws.Range("N27").Value = ""

# CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 23197.8
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 23197.8
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 69593.39999999999
$ws.Range("M81").Value = ""
$ws.Range("N81").Value = -71839.39999999999

# CUL row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 23197.8
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 23197.8
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 208780.2
$ws.Range("M84").Value = ""
$ws.Range("N84").Value = -220012.2

# CUL row 115
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 250000000
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = ""

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 738.8182
$ws.Range("J122").Value = 868.75
$ws.Range("L122").Value = 7818.75
$ws.Range("N122").Value = -12718.75

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 7485.091
$ws.Range("J137").Value = 34833
$ws.Range("L137").Value = 104499
$ws.Range("N137").Value = -114699

# GSM row 23
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 3000
$ws.Range("I23").Value = 3000
$ws.Range("K23").Value = 3000
$ws.Range("M23").Value = -2777

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 14130.75
$ws.Range("J46").Value = 16907.25
$ws.Range("L46").Value = 16907.25
$ws.Range("N46").Value = -17283.25

# LTW row 103
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 68999.5
$ws.Range("J103").Value = 68999.5
$ws.Range("L103").Value = 68999.5
$ws.Range("N103").Value = -71343.5

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4119.32
$ws.Range("I136").Value = 3690.1365
$ws.Range("J136").Value = 7266.6665
$ws.Range("K136").Value = 11070.4095
$ws.Range("L136").Value = 21799.9995
$ws.Range("M136").Value = -8520.4095
$ws.Range("N136").Value = -26899.9995

# WVR row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 16137.275
$ws.Range("I14").Value = 14293.177
$ws.Range("J14").Value = 18749.75
$ws.Range("K14").Value = 14293.177
$ws.Range("L14").Value = 18749.75
$ws.Range("M14").Value = -14125.177
$ws.Range("N14").Value = -19085.75

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1081.4667
$ws.Range("I100").Value = 1038.909
$ws.Range("K100").Value = 2077.818
$ws.Range("M100").Value = -1536.818

# WVR row 101
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 27602
$ws.Range("J101").Value = 27602
$ws.Range("L101").Value = 27602
$ws.Range("N101").Value = -34092

# WVR row 104
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 32249.75
$ws.Range("J104").Value = 32249.75
$ws.Range("L104").Value = 32249.75
$ws.Range("N104").Value = -39237.75

# WVR row 108
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 55612.5
$ws.Range("J108").Value = 55612.5
$ws.Range("L108").Value = 55612.5
$ws.Range("N108").Value = -63292.5

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 953.7273
$ws.Range("J113").Value = 779.5
$ws.Range("L113").Value = 2338.5
$ws.Range("N113").Value = -6678.5

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3248.4666
$ws.Range("I136").Value = 2209.7693
$ws.Range("K136").Value = 6629.3079
$ws.Range("M136").Value = -4079.3079
